# Results from July 22, 2020 05:56:00 PM America/Los_Angeles TZ run
#
# This run's GitHub API calls were rate-limited, so the "New York -- New
# York" row (row 4), which previously succeeded, now has no fetched data:
# the date/case/death/population-subset columns go back to being blank,
# the "includes Hispanic Black" flag (J4) reverts to False, and the status
# message reports the rate-limit error instead of "Success!". Separately,
# the Delaware row's (row 39) previous timeout error was replaced this run
# by a different failure (a numpy/pandas AttributeError) while fetching.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4:H4 and K4:L4 were populated numeric/text results; clear each back to
# an (empty-string) unfetched cell and strip the inherited formatting
# (B4 carried the date number format) so they match the other cleared
# cells in this column, e.g. K2/L2.
foreach ($addr in @("B4", "C4", "D4", "E4", "F4", "G4", "H4", "K4", "L4")) {
  $ws.Range($addr).Value = "'"
  $ws.Range($addr).ClearFormats()
}

# "Pct Includes Hispanic Black" flips back to False for this row.
$ws.Range("J4").Value = $false

# Status messages updated to this run's results.
$ws.Range("O4").Value = 'An error occurred. ... RateLimitExceededException(403, {''message'': "API rate limit exceeded for 132.145.200.60. (But here''s the good news: Authenticated requests get a higher rate limit. Check out the documentation for more details.)", ''documentation_url'': ''https://developer.github.com/v3/#rate-limiting''})'
$ws.Range("O39").Value = 'An error occurred. ... AttributeError("''numpy.float64'' object has no attribute ''split''")'
